# Apply "database sql update and interface information update" change.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 8: 医院介绍 (Hospital introduction) ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "医院介绍"
$ws.Range("C8").Value = "后台提供url，前端获取到医院详细信息和所有科室名字"
$ws.Range("E8").Value = "医院表hospital和科室表section"
$ws.Rows.Item(8).RowHeight = 28.8

# --- Row 9: 医生选择 (Doctor selection) ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "医生选择"
$ws.Range("C9").Value = "后台提供url，前端提供6中所选的科室名称， 后台根据科室名称返回前端该科室下所有医生列表"
$ws.Range("E9").Value = "医生表doctor"
$ws.Rows.Item(9).RowHeight = 43.2

# --- Row 10: 医生详细信息 (Doctor details) ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "医生详细信息"
$ws.Range("C10").Value = "前端提供7中所选的医生id，后端返回选择医生的所有详细信息"
$ws.Range("E10").Value = "医生表doctor"
$ws.Rows.Item(10).RowHeight = 28.8

# --- Row 11: 医院导航界面 (Hospital navigation interface) ---
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "医院导航界面"
$ws.Range("C11").Value = "后台提供医院交通官方网站交通"
$ws.Range("E11").Value = "医院表hospital"

# --- Update existing row 7 ("注册用户" / registration) ---
# Modify the interface-description text for the registration interface (column C, row 7)
# (done last so the new shared string is appended at the end, matching source order)
$ws.Range("C7").Value = "后台提供url，前端把注册信息传给后台，后台进行校验和数据插入，返回结果给前端（要进行校验！！！）"
$ws.Rows.Item(7).RowHeight = 57.6

# --- Sheet view: scroll position and selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C7").Select()
